$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- Fix the "CasesTab" query (cell B2): a stray cohort lookup/column was
# removed from the Cypher query text (OPTIONAL MATCH (co:cohort)... and the
# trailing `Cohort` return column), per "Fixed variables and query errors".
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Black and Tan Coonhound']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesQuery

# --- View state: the workbook was left scrolled/selected on B2 with the
# zoom bumped up to 130%.
$ws.Range("B2").Select() | Out-Null
$excel.ActiveWindow.Zoom = 130

# --- The row-2 wrap height shrinks because the query text lost two lines;
# rows 3/4 shift by a hair too (re-measured alongside it). Apply the
# resulting auto heights explicitly since this runtime doesn't recompute
# wrapped-text row heights on its own.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 216
